$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@('Hyundai Grand i10 Nios Magna','Rs. 5.83 Lakh','2026-02-16 11:06:31')
    ,@('Mercedes-Benz C-Class C 220 CDI Style','Rs. 14.14 Lakh','2026-02-16 11:06:31')
    ,@('Honda Jazz 1.2 V i VTEC','Rs. 4.21 Lakh','2026-02-16 11:06:31')
    ,@('Toyota Innova Crysta 2.8 ZX AT BSIV','Rs. 19.00 Lakh','2026-02-16 11:06:31')
    ,@('Hyundai i20 Sportz Plus Dual Tone BSIV','Rs. 5.78 Lakh','2026-02-16 11:06:31')
    ,@('Kia Sonet HTX Diesel iMT','Rs. 10.22 Lakh','2026-02-16 11:06:31')
    ,@('56 Used Hyundai i20','Rs. 69.50 Lakh','2026-02-16 11:06:31')
    ,@('53 Used Honda City','Rs. 4.15 Lakh','2026-02-16 11:06:31')
    ,@('45 Used Honda Amaze','Rs. 4.10 Lakh','2026-02-16 11:06:31')
    ,@('43 Used Maruti Swift','Rs. 3.95 Lakh','2026-02-16 11:06:31')
    ,@('40 Used Hyundai Creta','Rs. 3.90 Lakh','2026-02-16 11:06:31')
    ,@('39 Used Hyundai Grand i10','Rs. 3.70 Lakh','2026-02-16 11:06:31')
    ,@('33 Used Maruti Celerio','Rs. 4.85 Lakh','2026-02-16 11:06:31')
    ,@('33 Used Maruti Wagon R','Rs. 6.00 Lakh','2026-02-16 11:06:31')
    ,@('33 Used Kia Seltos','Rs. 3.15 Lakh','2026-02-16 11:06:31')
    ,@('Mercedes-Benz GLE 300d 4Matic','Rs. 9.00 Lakh','2026-02-16 11:06:31')
    ,@('Volkswagen Polo 1.2 MPI Highline','Rs. 9.50 Lakh','2026-02-16 11:06:31')
    ,@('Tata Tiago 1.2 Revotron XZA','Rs. 4.93 Lakh','2026-02-16 11:06:31')
    ,@('Hyundai Xcent 1.2 Kappa S AT','Rs. 7.50 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VDI BSIV','Rs. 10.75 Lakh','2026-02-16 11:06:31')
    ,@('Volkswagen Ameo 1.5 TDI Highline','Rs. 10.50 Lakh','2026-02-16 11:06:31')
    ,@('Honda City i VTEC S','Rs. 7.21 Lakh','2026-02-16 11:06:31')
    ,@('Honda City i VTEC CVT VX','Rs. 1.50 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VDI','Rs. 3.25 Lakh','2026-02-16 11:06:31')
    ,@('Mahindra XUV500 AT W10 FWD','Rs. 3.90 Lakh','2026-02-16 11:06:31')
    ,@('Honda City ZX MT','Rs. 3.65 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift ZXI AMT','Rs. 4.10 Lakh','2026-02-16 11:06:31')
    ,@('Mahindra XUV500 AT W10 FWD','Rs. 4.40 Lakh','2026-02-16 11:06:31')
    ,@('Mahindra XUV500 W9','Rs. 9.50 Lakh','2026-02-16 11:06:31')
    ,@('Honda City V MT','Rs. 7.45 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift ZXI','Rs. 4.55 Lakh','2026-02-16 11:06:31')
    ,@('Hyundai i10 Magna 1.1L','Rs. 4.75 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VXI','Rs. 4.52 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VDI BSIV','Rs. 4.75 Lakh','2026-02-16 11:06:31')
    ,@('Hyundai i10 Sportz 1.1L','Rs. 5.59 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift LXI Option','Rs. 3.10 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VXI','Rs. 3.15 Lakh','2026-02-16 11:06:31')
    ,@('Honda City ZX MT','Rs. 5.63 Lakh','2026-02-16 11:06:31')
    ,@('Honda City V CVT','Rs. 10.97 Lakh','2026-02-16 11:06:31')
    ,@('Honda City 1.5 V MT Exclusive','Rs. 8.75 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift RS VDI','Rs. 4.95 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift ZXI','Rs. 6.53 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VXI Optional','Rs. 11.06 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VDI','Rs. 6.80 Lakh','2026-02-16 11:06:31')
    ,@('Hyundai i10 Magna 1.1','Rs. 3.25 Lakh','2026-02-16 11:06:31')
    ,@('Hyundai i10 Magna 1.1','Rs. 3.75 Lakh','2026-02-16 11:06:31')
    ,@('Honda City VX MT','Rs. 12.09 Lakh','2026-02-16 11:06:31')
    ,@('Honda City V','Rs. 10.50 Lakh','2026-02-16 11:06:31')
    ,@('Honda City V 2 Airbag','Rs. 11.39 Lakh','2026-02-16 11:06:31')
    ,@('Honda City i VTEC CVT SV','Rs. 3.49 Lakh','2026-02-16 11:06:31')
    ,@('Honda City i-VTEC V','Rs. 4.90 Lakh','2026-02-16 11:06:31')
    ,@('Honda City VX MT','Rs. 3.62 Lakh','2026-02-16 11:06:31')
    ,@('Honda City i VTEC VX','Rs. 4.77 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VXI','Rs. 3.37 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VDI','Rs. 4.14 Lakh','2026-02-16 11:06:31')
    ,@('Honda City ZX MT','Rs. 4.80 Lakh','2026-02-16 11:06:31')
    ,@('Honda City VX MT','Rs. 4.25 Lakh','2026-02-16 11:06:31')
    ,@('Honda City ZX MT','Rs. 5.80 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VXI','Rs. 6.93 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift ZDI BSIV','Rs. 6.95 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift ZXI','Rs. 3.15 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift AMT ZDI','Rs. 13.95 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift LXI Optional-O','Rs. 3.65 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VDI BSIV','Rs. 3.45 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VXI AMT','Rs. 27.00 Lakh','2026-02-16 11:06:31')
    ,@('Honda City V MT','Rs. 9.50 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift AMT VXI','Rs. 16.75 Lakh','2026-02-16 11:06:31')
    ,@('Honda City i-VTEC CVT ZX','Rs. 10.75 Lakh','2026-02-16 11:06:31')
    ,@('Honda City i-VTEC V','Rs. 2.80 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VDI','Rs. 2.95 Lakh','2026-02-16 11:06:31')
    ,@('Toyota Innova 2.5 G (Diesel) 7 Seater BS IV','Rs. 16.50 Lakh','2026-02-16 11:06:31')
    ,@('Maruti Swift VDI','Rs. 9.75 Lakh','2026-02-16 11:06:31')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}